# Update TPM-derived NATMI metrics with newly recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.385511999999999
$ws.Range("H2").Value = 16.156536
$ws.Range("I2").Value = 0.02736372477514656
$ws.Range("J2").Value = 0.02736372477514657
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01506066666666667
$ws.Range("N2").Value = 0.045182
$ws.Range("O2").Value = 0.2394078124668829
$ws.Range("P2").Value = 0.2394078124668829
$ws.Range("Q2").Value = 0.08110940106133332
$ws.Range("R2").Value = 0.729984609552
$ws.Range("S2").Value = 0.006551089489363685
$ws.Range("T2").Value = 0.006551089489363686

# Row 3
$ws.Range("G3").Value = 5.385511999999999
$ws.Range("H3").Value = 16.156536
$ws.Range("I3").Value = 0.02736372477514656
$ws.Range("J3").Value = 0.02736372477514657
$ws.Range("O3").Value = 0.7605921875331172
$ws.Range("P3").Value = 0.7605921875331172
$ws.Range("Q3").Value = 0.2576823878346666
$ws.Range("R3").Value = 2.319141490512
$ws.Range("S3").Value = 0.02081263528578288
$ws.Range("T3").Value = 0.02081263528578288

# Row 4
$ws.Range("I4").Value = 0.7812411799860843
$ws.Range("J4").Value = 0.7812411799860843
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01506066666666667
$ws.Range("N4").Value = 0.045182
$ws.Range("O4").Value = 0.2394078124668829
$ws.Range("P4").Value = 0.2394078124668829
$ws.Range("Q4").Value = 2.315693667942222
$ws.Range("R4").Value = 20.84124301148
$ws.Range("S4").Value = 0.1870352419095148
$ws.Range("T4").Value = 0.1870352419095148

# Row 5
$ws.Range("I5").Value = 0.7812411799860843
$ws.Range("J5").Value = 0.7812411799860843
$ws.Range("O5").Value = 0.7605921875331172
$ws.Range("P5").Value = 0.7605921875331172
$ws.Range("S5").Value = 0.5942059380765696
$ws.Range("T5").Value = 0.5942059380765696

# Row 6
$ws.Range("I6").Value = 0.1913950952387691
$ws.Range("J6").Value = 0.1913950952387691
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01506066666666667
$ws.Range("N6").Value = 0.045182
$ws.Range("O6").Value = 0.2394078124668829
$ws.Range("P6").Value = 0.2394078124668829
$ws.Range("Q6").Value = 0.5673182897597777
$ws.Range("R6").Value = 5.105864607838
$ws.Range("S6").Value = 0.04582148106800442
$ws.Range("T6").Value = 0.04582148106800441

# Row 7
$ws.Range("I7").Value = 0.1913950952387691
$ws.Range("J7").Value = 0.1913950952387691
$ws.Range("O7").Value = 0.7605921875331172
$ws.Range("P7").Value = 0.7605921875331172
$ws.Range("S7").Value = 0.1455736141707647
$ws.Range("T7").Value = 0.1455736141707647
